# Generate Report for handback
# Updates the localization status workbook: zh-cn and de-de sheets get a
# second "handback" round recorded (Latest Target File / Latest Handback
# File / Latest Handback DateTime), and the Status column moves from
# "Ready for handoff" to "Handed back: in sync with en-us" everywhere
# (Overview sheet picks this up automatically because it shares the same
# text).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-us"

# ---------------------------------------------------------------------
# Overview sheet: status column text changes for both language rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("E2").Value = "3471dafb-1ade-4f2f-9627-f17e7208feec.md"
$wsZh.Range("F2").Value = "3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-01-08 19:57:58"

$wsZh.Range("B3").Value = $statusText
$wsZh.Range("E3").Value = "980d8046-22b9-43b4-98c5-ca34348e8d26.md"
$wsZh.Range("F3").Value = "980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-01-08 19:57:58"

# Hyperlinks for the newly-populated cells (mirrors columns A / C).
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/aff2c7e1e52b231bccfbc9703e1d7c07c781babd/e2e/3471dafb-1ade-4f2f-9627-f17e7208feec.md", [Type]::Missing, [Type]::Missing, "3471dafb-1ade-4f2f-9627-f17e7208feec.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d23a314eea8c0e952da67532ddce34b792c041a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/aff2c7e1e52b231bccfbc9703e1d7c07c781babd/e2e/980d8046-22b9-43b4-98c5-ca34348e8d26.md", [Type]::Missing, [Type]::Missing, "980d8046-22b9-43b4-98c5-ca34348e8d26.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d23a314eea8c0e952da67532ddce34b792c041a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("E2").Value = "3471dafb-1ade-4f2f-9627-f17e7208feec.md"
$wsDe.Range("F2").Value = "3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.de-de.xlf"
$wsDe.Range("G2").Value = "2016-01-08 19:58:16"

$wsDe.Range("B3").Value = $statusText
$wsDe.Range("E3").Value = "980d8046-22b9-43b4-98c5-ca34348e8d26.md"
$wsDe.Range("F3").Value = "980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.de-de.xlf"
$wsDe.Range("G3").Value = "2016-01-08 19:58:16"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/aff2c7e1e52b231bccfbc9703e1d7c07c781babd/e2e/3471dafb-1ade-4f2f-9627-f17e7208feec.md", [Type]::Missing, [Type]::Missing, "3471dafb-1ade-4f2f-9627-f17e7208feec.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a6dbe3c34f140f4d4b489f2421bbc398c7c296d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.de-de.xlf", [Type]::Missing, [Type]::Missing, "3471dafb-1ade-4f2f-9627-f17e7208feec.f65efd1ce3074a32632c01c99f2aeb1d7a855be0.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/aff2c7e1e52b231bccfbc9703e1d7c07c781babd/e2e/980d8046-22b9-43b4-98c5-ca34348e8d26.md", [Type]::Missing, [Type]::Missing, "980d8046-22b9-43b4-98c5-ca34348e8d26.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a6dbe3c34f140f4d4b489f2421bbc398c7c296d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.de-de.xlf", [Type]::Missing, [Type]::Missing, "980d8046-22b9-43b4-98c5-ca34348e8d26.af298d2d4527b3556d684643f18890ed97c77661.de-de.xlf") | Out-Null

# Give the new cells the same "HyperLink" look as the existing link cells.
$wsZh.Range("E2:F3").Style = "HyperLink"
$wsDe.Range("E2:F3").Style = "HyperLink"
